# Update the "cryptos" price list snapshot (Coin / Link / Price / Volume(1h) columns).
# Cells in column D that look like plain numbers are forced to remain text (format "@")
# because the sheet stores prices such as "43.695.22" or "94.94" as literal strings,
# and Excel would otherwise silently reinterpret/round them as numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.695.22'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '2.299.76'
$ws.Range("E3").Value = '  +2.10%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '94.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.05%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.610'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0938'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.81'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.34%  '
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").Value = '2.645.27'
$ws.Range("E14").Value = '  +9.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.856'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.63%  '
$ws.Range("D17").Value = '2.301.28'
$ws.Range("E17").Value = '  +2.22%  '
$ws.Range("D18").Value = '43.673.92'
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000108'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.97'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E26").Value = '  -1.23%  '
$ws.Range("E27").Value = '  -3.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.51%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0888'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.98%  '
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("E36").Value = '  -5.26%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0354'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.43'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.26'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -13.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.37'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.238'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.58%  '
$ws.Range("E42").Value = '  +17.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.68%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.102'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '99.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.97%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.53'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.09%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.522.88'
$ws.Range("E51").Value = '  +2.97%  '
